# Automatische test-sync: 2025-06-23 18:18:50
# Adds a new incoming mail-log row (row 8) to the "Logs" sheet and
# updates the "Dashboard" category-count table to reflect the new totals.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- 1. Append the new mail-log entry on row 8 ------------------------
$logs.Range("A8").Value = "Korting voor wederverkopers?"
$logs.Range("B8").Value = "mailmind.test@zohomail.eu"
$logs.Range("C8").Value = "Biedt u speciale prijzen voor wederverkopers?"
$logs.Range("D8").Value = "Offerte / Prijsaanvraag"
$logs.Range("F8").Value = "2025-06-23 18:18:07"
$logs.Range("G8").Value = "Nee"

# --- 2. Refresh the Dashboard category/count table ---------------------
# "Offerte / Prijsaanvraag" now has 2 occurrences, so it moves up to row 3
# (right after "IT / Technisch probleem", which stays at 2), and the
# remaining categories shift down one row, keeping their counts of 1.
$dashboard.Range("A3").Value = "Offerte / Prijsaanvraag"
$dashboard.Range("B3").Value = 2

$dashboard.Range("A4").Value = "Bestelling / Levering"
$dashboard.Range("B4").Value = 1

$dashboard.Range("A5").Value = "Sollicitatie / Vacature"
$dashboard.Range("B5").Value = 1

$dashboard.Range("A6").Value = "Factuur / Administratie"
$dashboard.Range("B6").Value = 1

# --- 3. Extend the conditional-formatting ranges to include the new row -
$catFcs = $logs.Range("D2:D7").FormatConditions
for ($i = 1; $i -le $catFcs.Count; $i++) {
    $catFcs.Item($i).ModifyAppliesToRange($logs.Range("D2:D8"))
}

$answeredFcs = $logs.Range("G2:G7").FormatConditions
for ($i = 1; $i -le $answeredFcs.Count; $i++) {
    $answeredFcs.Item($i).ModifyAppliesToRange($logs.Range("G2:G8"))
}
